# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-04-02 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-03 Thursday", 2) | Out-Null

# Update the multiplication-fact table. Addressed by (row, column) since
# several cells share identical source text, so a blanket Find/Replace
# would be ambiguous.
$table = $d.Tables(1)

$updates = @(
    @{Row=1;  Col=1; Text="880×6=5280"},
    @{Row=1;  Col=2; Text="727×9=6543"},
    @{Row=1;  Col=3; Text="771×5=3855"},
    @{Row=1;  Col=4; Text="418×8=3344"},
    @{Row=1;  Col=5; Text="415×5=2075"},

    @{Row=5;  Col=1; Text="641×5=3205"},
    @{Row=5;  Col=2; Text="581×2=1162"},
    @{Row=5;  Col=3; Text="106×3=318"},
    @{Row=5;  Col=4; Text="521×4=2084"},
    @{Row=5;  Col=5; Text="835×9=7515"},

    @{Row=10; Col=1; Text="690×8=5520"},
    @{Row=10; Col=2; Text="147×2=294"},
    @{Row=10; Col=3; Text="535×4=2140"},
    @{Row=10; Col=4; Text="194×2=388"},
    @{Row=10; Col=5; Text="405×5=2025"},

    @{Row=15; Col=1; Text="597×7=4179"},
    @{Row=15; Col=2; Text="920×7=6440"},
    @{Row=15; Col=3; Text="770×5=3850"},
    @{Row=15; Col=4; Text="656×9=5904"},
    @{Row=15; Col=5; Text="314×6=1884"},

    @{Row=20; Col=1; Text="804×7=5628"},
    @{Row=20; Col=2; Text="309×5=1545"},
    @{Row=20; Col=3; Text="355×8=2840"},
    @{Row=20; Col=4; Text="697×9=6273"},
    @{Row=20; Col=5; Text="216×9=1944"}
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
